$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ----
$ws.Range("E2").Value = 24.60000000000041
$ws.Range("H2").Value = [double]"2.508978586723518e-16"
$ws.Range("I2").Value = 0.657991530246405
$ws.Range("K2").Value = 51.81791299112101
$ws.Range("L2").Value = "[41.58485233331823, 62.05097364892378]"
$ws.Range("O2").Value = 1.691868716347656
$ws.Range("P2").Value = "[1.46544762419704, 1.9182898084982725]"
$ws.Range("S2").Value = 64.18536413143424
$ws.Range("T2").Value = "[57.68538788768399, 70.6853403751845]"
$ws.Range("W2").Value = 17.97597597597627
$ws.Range("X2").Value = 17.08948948948977
$ws.Range("Y2").Value = 18.86246246246278

# ---- Row 3 ----
$ws.Range("B3").Value = 1
$ws.Range("E3").Value = 23.87000000000029
$ws.Range("G3").Value = [double]"1.110223024625157e-16"
$ws.Range("H3").Value = [double]"2.508978586723518e-16"
$ws.Range("I3").ClearContents()
$ws.Range("K3").Value = 50.18057517519721
$ws.Range("L3").Value = "[36.17393531694384, 64.18721503345057]"
$ws.Range("M3").Value = [double]"1.854050246663519e-11"
$ws.Range("N3").Value = [double]"1.854050246663519e-11"
$ws.Range("O3").Value = 0.7107106503616549
$ws.Range("P3").Value = "[0.4339737599553466, 0.9874475407679633]"
$ws.Range("Q3").Value = [double]"8.432999072383041e-07"
$ws.Range("R3").Value = [double]"8.432999072383041e-07"
$ws.Range("S3").Value = 65.59253413924077
$ws.Range("T3").Value = "[58.16686450991066, 73.01820376857088]"
$ws.Range("W3").Value = 21.16998998999025
$ws.Range("X3").Value = 20.11865865865891
$ws.Range("Y3").Value = 22.22132132132159
